$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 104449393
$ws.Range("B2").Value = 78570
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("K2").ClearContents() | Out-Null
$ws.Range("L2").ClearContents() | Out-Null
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").ClearContents() | Out-Null
$ws.Range("Q2").Value = 418188.408122587
$ws.Range("R2").Value = 7018072.943679515
$ws.Range("AC2").ClearContents() | Out-Null

# Row 3
$ws.Range("A3").Value = 104449298
$ws.Range("Q3").Value = 418152.1433075544
$ws.Range("R3").Value = 7018755.833866266

# Row 4
$ws.Range("A4").Value = 104449306
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 418116.3208070688
$ws.Range("R4").Value = 7018906.624424814

# Row 5
$ws.Range("A5").Value = 104449305
$ws.Range("M5").Value = "färska spår"
$ws.Range("Q5").Value = 418113.1107625436
$ws.Range("R5").Value = 7018904.455793464

# Row 6
$ws.Range("A6").Value = 104449333
$ws.Range("Q6").Value = 418224.807980529
$ws.Range("R6").Value = 7018298.906277624

# Row 7
$ws.Range("A7").Value = 104449307
$ws.Range("M7").Value = "äldre spår"
$ws.Range("Q7").Value = 418106.0762497109
$ws.Range("R7").Value = 7018911.38607322
$ws.Range("AC7").Value = "ringhack"

# Row 8
$ws.Range("A8").Value = 104449297
$ws.Range("Q8").Value = 418163.1633477406
$ws.Range("R8").Value = 7018746.101364438

# Row 9
$ws.Range("A9").Value = 104449308
$ws.Range("B9").Value = 56395
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("K9").ClearContents() | Out-Null
$ws.Range("L9").ClearContents() | Out-Null
$ws.Range("M9").Value = "färska spår"
$ws.Range("N9").ClearContents() | Out-Null
$ws.Range("Q9").Value = 418207.1051796933
$ws.Range("R9").Value = 7019144.644948276
$ws.Range("AC9").Value = "Påbörjat bo?"
